$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.356.79"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").Value = "2.079.14"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.97"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.663"
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.97"
$ws.Range("E8").Value = "  +25.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "62.39"
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.385"
$ws.Range("E10").Value = "  +5.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0752"
$ws.Range("E11").Value = "  +2.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.106"
$ws.Range("E12").Value = "  +7.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.46"
$ws.Range("E13").Value = "  +5.53%  "
$ws.Range("D14").Value = "2.389.84"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.846"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("E16").Value = "  +4.57%  "
$ws.Range("D17").Value = "2.086.39"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").Value = "37.271.45"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.18"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.76"
$ws.Range("E20").Value = "  +14.58%  "
$ws.Range("D21").Value = "0.0₃0849"
$ws.Range("E21").Value = "  +3.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "240.65"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("E23").Value = "  +4.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.47"
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.42"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.10"
$ws.Range("E27").Value = "  +2.74%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.22"
$ws.Range("E28").Value = "  +2.70%  "
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("B31").Value = "Gas"
$ws.Range("C31").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.07"
$ws.Range("E31").Value = "  +3.29%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.10"
$ws.Range("E32").Value = "  +21.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.57"
$ws.Range("E33").Value = "  +2.94%  "
$ws.Range("E34").Value = "  +5.35%  "
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.31"
$ws.Range("E36").Value = "  +5.84%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.86"
$ws.Range("E38").Value = "  -1.43%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.28"
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("E40").Value = "  -1.03%  "
$ws.Range("E41").Value = "  +4.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0995"
$ws.Range("E42").Value = "  +19.14%  "
$ws.Range("E43").Value = "  +9.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.17"
$ws.Range("E44").Value = "  -1.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.89"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("E46").Value = "  +111.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.80"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").Value = "1.326.91"
$ws.Range("E48").Value = "  -2.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.94"
$ws.Range("E49").Value = "  +3.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.36"
$ws.Range("E50").Value = "  +4.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.03"
$ws.Range("E51").Value = "  +9.31%  "
